$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A holds the "Beteckning" key for every data row; use it to find the
# last populated data row (data starts on row 2, row 1 is the header).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 3).Value = 45172
}
